# Generate Report for Handoff
#
# The "c05037ce-57fe-4517-b9e1-e8532af8799d" item has finished translation and
# is now "Ready for handoff" (with a fresh handback timestamp), while the
# "45adacbf-760e-419c-8f27-a66d66377ffe" item - which used to occupy row 7 with
# that same status - moves up into row 6 carrying its original (older)
# timestamps. This script rewrites rows 6 and 7 on the Overview, zh-cn and
# de-de worksheets to reflect the new row order / values, and fixes up the
# corresponding hyperlink display text.

$wb = $excel.ActiveWorkbook

function Set-HyperlinkDisplay {
    param(
        $ws,
        [string]$address,
        [string]$text
    )
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $address) {
            $hl.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 6 now holds the 45adacbf item (status: Ready for handoff)
$wsOverview.Range("A6").Value = "45adacbf-760e-419c-8f27-a66d66377ffe.md"
$wsOverview.Range("B6").Value = "e2e\45adacbf-760e-419c-8f27-a66d66377ffe.md"
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "2017-02-17 07:59:38"

# Row 7 now holds the c05037ce item (status: Ready for handoff, new timestamp)
$wsOverview.Range("A7").Value = "c05037ce-57fe-4517-b9e1-e8532af8799d.md"
$wsOverview.Range("B7").Value = "e2e\c05037ce-57fe-4517-b9e1-e8532af8799d.md"
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = "2017-02-17 08:06:20"

Set-HyperlinkDisplay $wsOverview '$B$6' "e2e\45adacbf-760e-419c-8f27-a66d66377ffe.md"
Set-HyperlinkDisplay $wsOverview '$B$7' "e2e\c05037ce-57fe-4517-b9e1-e8532af8799d.md"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A6").Value = "45adacbf-760e-419c-8f27-a66d66377ffe.md"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("G6").Value = "45adacbf-760e-419c-8f27-a66d66377ffe.6e3d7dc9fc5e0000dbe65cfe718c97203585a820.zh-cn.xlf"
$wsZhCn.Range("H6").Value = "2017-02-17 07:59:22"

$wsZhCn.Range("A7").Value = "c05037ce-57fe-4517-b9e1-e8532af8799d.md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("G7").Value = "c05037ce-57fe-4517-b9e1-e8532af8799d.65d45af1316fd47f5e4d4ea03c88ebac9c0dc445.zh-cn.xlf"
$wsZhCn.Range("H7").Value = "2017-02-17 08:06:03"

Set-HyperlinkDisplay $wsZhCn '$A$6' "45adacbf-760e-419c-8f27-a66d66377ffe.md"
Set-HyperlinkDisplay $wsZhCn '$A$7' "c05037ce-57fe-4517-b9e1-e8532af8799d.md"

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A6").Value = "45adacbf-760e-419c-8f27-a66d66377ffe.md"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("G6").Value = "45adacbf-760e-419c-8f27-a66d66377ffe.6e3d7dc9fc5e0000dbe65cfe718c97203585a820.de-de.xlf"
$wsDeDe.Range("H6").Value = "2017-02-17 07:59:38"

$wsDeDe.Range("A7").Value = "c05037ce-57fe-4517-b9e1-e8532af8799d.md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("G7").Value = "c05037ce-57fe-4517-b9e1-e8532af8799d.65d45af1316fd47f5e4d4ea03c88ebac9c0dc445.de-de.xlf"
$wsDeDe.Range("H7").Value = "2017-02-17 08:06:20"

Set-HyperlinkDisplay $wsDeDe '$A$6' "45adacbf-760e-419c-8f27-a66d66377ffe.md"
Set-HyperlinkDisplay $wsDeDe '$A$7' "c05037ce-57fe-4517-b9e1-e8532af8799d.md"

$wb.Save()
